# Aggiornamento fino a 27/05: append 14 new daily rows (256-269) below the
# existing data table, re-using the formatting of the template row (255) for
# column A (date cell border/alignment/bold/date-time format) while columns
# B-D keep the default (unstyled) look already used throughout the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateRow = "A255:D255"

$ws.Range($templateRow).Copy($ws.Range("A256:D256"))
$ws.Cells.Item(256, 1).Value = 44330
$ws.Cells.Item(256, 2).Value = 2
$ws.Cells.Item(256, 3).Value = 16
$ws.Cells.Item(256, 4).Value = 66.32949175026947

$ws.Range($templateRow).Copy($ws.Range("A257:D257"))
$ws.Cells.Item(257, 1).Value = 44331
$ws.Cells.Item(257, 2).Value = 1
$ws.Cells.Item(257, 3).Value = 12
$ws.Cells.Item(257, 4).Value = 49.7471188127021

$ws.Range($templateRow).Copy($ws.Range("A258:D258"))
$ws.Cells.Item(258, 1).Value = 44332
$ws.Cells.Item(258, 2).Value = 0
$ws.Cells.Item(258, 3).Value = 7
$ws.Cells.Item(258, 4).Value = 29.01915264074289

$ws.Range($templateRow).Copy($ws.Range("A259:D259"))
$ws.Cells.Item(259, 1).Value = 44333
$ws.Cells.Item(259, 2).Value = 2
$ws.Cells.Item(259, 3).Value = 9
$ws.Cells.Item(259, 4).Value = 37.31033910952657

$ws.Range($templateRow).Copy($ws.Range("A260:D260"))
$ws.Cells.Item(260, 1).Value = 44334
$ws.Cells.Item(260, 2).Value = 2
$ws.Cells.Item(260, 3).Value = 10
$ws.Cells.Item(260, 4).Value = 41.45593234391841

$ws.Range($templateRow).Copy($ws.Range("A261:D261"))
$ws.Cells.Item(261, 1).Value = 44335
$ws.Cells.Item(261, 2).Value = 1
$ws.Cells.Item(261, 3).Value = 10
$ws.Cells.Item(261, 4).Value = 41.45593234391841

$ws.Range($templateRow).Copy($ws.Range("A262:D262"))
$ws.Cells.Item(262, 1).Value = 44336
$ws.Cells.Item(262, 2).Value = 0
$ws.Cells.Item(262, 3).Value = 8
$ws.Cells.Item(262, 4).Value = 33.16474587513473

$ws.Range($templateRow).Copy($ws.Range("A263:D263"))
$ws.Cells.Item(263, 1).Value = 44337
$ws.Cells.Item(263, 2).Value = 1
$ws.Cells.Item(263, 3).Value = 7
$ws.Cells.Item(263, 4).Value = 29.01915264074289

$ws.Range($templateRow).Copy($ws.Range("A264:D264"))
$ws.Cells.Item(264, 1).Value = 44338
$ws.Cells.Item(264, 2).Value = 1
$ws.Cells.Item(264, 3).Value = 7
$ws.Cells.Item(264, 4).Value = 29.01915264074289

$ws.Range($templateRow).Copy($ws.Range("A265:D265"))
$ws.Cells.Item(265, 1).Value = 44339
$ws.Cells.Item(265, 2).Value = 0
$ws.Cells.Item(265, 3).Value = 7
$ws.Cells.Item(265, 4).Value = 29.01915264074289

$ws.Range($templateRow).Copy($ws.Range("A266:D266"))
$ws.Cells.Item(266, 1).Value = 44340
$ws.Cells.Item(266, 2).Value = 4
$ws.Cells.Item(266, 3).Value = 9
$ws.Cells.Item(266, 4).Value = 37.31033910952657

$ws.Range($templateRow).Copy($ws.Range("A267:D267"))
$ws.Cells.Item(267, 1).Value = 44341
$ws.Cells.Item(267, 2).Value = 3
$ws.Cells.Item(267, 3).Value = 10
$ws.Cells.Item(267, 4).Value = 41.45593234391841

$ws.Range($templateRow).Copy($ws.Range("A268:D268"))
$ws.Cells.Item(268, 1).Value = 44342
$ws.Cells.Item(268, 2).Value = 0
$ws.Cells.Item(268, 3).Value = 9
$ws.Cells.Item(268, 4).Value = 37.31033910952657

$ws.Range($templateRow).Copy($ws.Range("A269:D269"))
$ws.Cells.Item(269, 1).Value = 44343
$ws.Cells.Item(269, 2).Value = 1
$ws.Cells.Item(269, 3).Value = 10
$ws.Cells.Item(269, 4).Value = 41.45593234391841
